$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (e.g. "1.003", "24.794.74") that must stay
# text, not be auto-coerced to numbers - so force text format first.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.794.74"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.663.87"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.53"
$ws.Range("E5").Value = "  +7.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3643"
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.30"
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3256"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.134"
$ws.Range("E10").Value = "  +2.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07073"
$ws.Range("E11").Value = "  +2.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.071"
$ws.Range("E13").Value = "  +2.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.43"
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.663.41"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.592"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001048"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06646"
$ws.Range("E18").Value = "  +2.00%  "
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.42"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.919"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.75"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.52"
$ws.Range("E23").Value = "  +5.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.793.15"
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.431"
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.54"
$ws.Range("E27").Value = "  +3.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.65"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.846.42"
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.15"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.169"
$ws.Range("E31").Value = "  +7.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.063"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.678"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08487"
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.647"
$ws.Range("E35").Value = "  -1.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.15"
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06229"
$ws.Range("E37").Value = "  +4.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.157"
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("E39").Value = "  +3.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.236"
$ws.Range("E40").Value = "  +3.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2082"
$ws.Range("E41").Value = "  +2.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.216"
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5921"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.47"
$ws.Range("E45").Value = "  +8.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.838"
$ws.Range("E46").Value = "  +3.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5645"
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.57"
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.945"
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06969"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.192"
$ws.Range("E51").Value = "  +4.43%  "
